$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("SwateTemplateMetadata")

# Fix the Table name on the SwateTemplateMetadata sheet (row 6 "Table" -> B6)
$ws2.Range("B6").Value = "annotationTablePrettyGecko11"

# Update the active-cell selection on Sheet1 to match the saved cursor position
$ws1.Select()
$ws1.Range("B2").Select()

# Update the active-cell selection on SwateTemplateMetadata and leave it as the
# active (tab-selected) sheet, matching the original workbook state
$ws2.Select()
$ws2.Range("B7").Select()
